# Auto-generated Excel COM-interop script
# Applies updated market-price snapshot values to the Tonberry Profits workbook
# (data pulled from an external price API on a scheduled run; no formulas involved)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 191.8421
$ws.Range("I33").Value = 230.83333
$ws.Range("K33").Value = 230.83333
$ws.Range("M33").Value = -1.833329999999989
$ws.Range("H64").Value = 3666.3333
$ws.Range("J64").Value = 3999.5
$ws.Range("L64").Value = 3999.5
$ws.Range("N64").Value = -4495.5
$ws.Range("H67").Value = 3666.3333
$ws.Range("J67").Value = 3999.5
$ws.Range("L67").Value = 3999.5
$ws.Range("N67").Value = -5715.5
$ws.Range("H98").Value = 1386
$ws.Range("I98").Value = 1251.5
$ws.Range("K98").Value = 1251.5
$ws.Range("M98").Value = 246.5
$ws.Range("H100").Value = 1756.2858
$ws.Range("J100").Value = 2999.5
$ws.Range("L100").Value = 2999.5
$ws.Range("N100").Value = -4081.5
$ws.Range("H113").Value = 24166
$ws.Range("I113").Value = 30685.572
$ws.Range("K113").Value = 30685.572
$ws.Range("M113").Value = -27431.572
$ws.Range("H121").Value = 1229.6
$ws.Range("J121").Value = 1229.6
$ws.Range("L121").Value = 3688.8
$ws.Range("N121").Value = -7182.799999999999
$ws.Range("H122").Value = 1386
$ws.Range("I122").Value = 1251.5
$ws.Range("K122").Value = 3754.5
$ws.Range("M122").Value = -1304.5
$ws.Range("H125").Value = 1196
$ws.Range("I125").Value = 1900
$ws.Range("K125").Value = 17100
$ws.Range("M125").Value = -14640

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3830.8696
$ws.Range("I32").Value = 2421.182
$ws.Range("J32").Value = 7409.3076
$ws.Range("K32").Value = 2421.182
$ws.Range("L32").Value = 7409.3076
$ws.Range("M32").Value = -2134.182
$ws.Range("N32").Value = -7983.3076
$ws.Range("H38").Value = 3000
$ws.Range("I38").Value = 3000
$ws.Range("K38").Value = 3000
$ws.Range("M38").Value = -2533
$ws.Range("H45").Value = 2326.25
$ws.Range("I45").Value = 736
$ws.Range("K45").Value = 736
$ws.Range("M45").Value = -359
$ws.Range("H46").Value = 3247
$ws.Range("I46").Value = 3358.75
$ws.Range("J46").Value = 2800
$ws.Range("K46").Value = 3358.75
$ws.Range("L46").Value = 2800
$ws.Range("M46").Value = -3039.75
$ws.Range("N46").Value = -3438
$ws.Range("H61").Value = 7431
$ws.Range("I61").Value = 4965.3335
$ws.Range("J61").Value = 9896.666999999999
$ws.Range("K61").Value = 4965.3335
$ws.Range("L61").Value = 9896.666999999999
$ws.Range("M61").Value = -4753.3335
$ws.Range("N61").Value = -10320.667
$ws.Range("H63").Value = 8997
$ws.Range("J63").Value = 7995
$ws.Range("L63").Value = 7995
$ws.Range("N63").Value = -9367
$ws.Range("H66").Value = 8997
$ws.Range("J66").Value = 7995
$ws.Range("L66").Value = 39975
$ws.Range("N66").Value = -46839
$ws.Range("H88").Value = 2601.2856
$ws.Range("I88").Value = 1824.3334
$ws.Range("J88").Value = 3999.8
$ws.Range("K88").Value = 1824.3334
$ws.Range("L88").Value = 3999.8
$ws.Range("M88").Value = -1418.3334
$ws.Range("N88").Value = -4811.8
$ws.Range("H91").Value = 2601.2856
$ws.Range("I91").Value = 1824.3334
$ws.Range("J91").Value = 3999.8
$ws.Range("K91").Value = 1824.3334
$ws.Range("L91").Value = 3999.8
$ws.Range("M91").Value = -420.3334
$ws.Range("N91").Value = -6807.8
$ws.Range("H97").Value = 593.86664
$ws.Range("I97").Value = 625.9167
$ws.Range("K97").Value = 625.9167
$ws.Range("M97").Value = -129.9167
$ws.Range("H110").Value = 3166.6667
$ws.Range("H122").Value = 2090.0908
$ws.Range("I122").Value = 2134.5557
$ws.Range("K122").Value = 6403.6671
$ws.Range("M122").Value = -3953.6671
$ws.Range("H136").Value = 7431
$ws.Range("I136").Value = 4965.3335
$ws.Range("J136").Value = 9896.666999999999
$ws.Range("K136").Value = 14896.0005
$ws.Range("L136").Value = 29690.001
$ws.Range("M136").Value = -12346.0005
$ws.Range("N136").Value = -34790.001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 13280
$ws.Range("J61").Value = 13280
$ws.Range("L61").Value = 13280
$ws.Range("N61").Value = -13906
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H134").Value = 6512.4326
$ws.Range("I134").Value = 6931.7144
$ws.Range("K134").Value = 20795.1432
$ws.Range("M134").Value = -18260.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2557.2856
$ws.Range("I31").Value = 1021.26086
$ws.Range("J31").Value = 5501.3335
$ws.Range("K31").Value = 1021.26086
$ws.Range("L31").Value = 5501.3335
$ws.Range("M31").Value = -726.26086
$ws.Range("N31").Value = -6091.3335
$ws.Range("H34").Value = 2557.2856
$ws.Range("I34").Value = 1021.26086
$ws.Range("J34").Value = 5501.3335
$ws.Range("K34").Value = 1021.26086
$ws.Range("L34").Value = 5501.3335
$ws.Range("M34").Value = -819.26086
$ws.Range("N34").Value = -5905.3335
$ws.Range("H58").Value = 1469.75
$ws.Range("J58").Value = 3000
$ws.Range("L58").Value = 3000
$ws.Range("N58").Value = -3406
$ws.Range("H93").Value = 16199.75
$ws.Range("J93").Value = 50000
$ws.Range("L93").Value = 50000
$ws.Range("N93").Value = -53744
$ws.Range("H99").Value = 2218.8
$ws.Range("I99").Value = 1549.5
$ws.Range("J99").Value = 2665
$ws.Range("K99").Value = 1549.5
$ws.Range("L99").Value = 2665
$ws.Range("M99").Value = -51.5
$ws.Range("N99").Value = -5661
$ws.Range("H126").Value = 2218.8
$ws.Range("I126").Value = 1549.5
$ws.Range("J126").Value = 2665
$ws.Range("K126").Value = 4648.5
$ws.Range("L126").Value = 7995
$ws.Range("M126").Value = -2178.5
$ws.Range("N126").Value = -12935
$ws.Range("H136").Value = 1469.75
$ws.Range("J136").Value = 3000
$ws.Range("L136").Value = 9000
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 110.5
$ws.Range("I23").Value = 40
$ws.Range("K23").Value = 120
$ws.Range("M23").Value = 115
$ws.Range("H88").Value = 4999.6665
$ws.Range("I88").Value = 3000
$ws.Range("J88").Value = 5999.5
$ws.Range("K88").Value = 9000
$ws.Range("L88").Value = 17998.5
$ws.Range("M88").Value = -8572
$ws.Range("N88").Value = -18854.5
$ws.Range("H91").Value = 4999.6665
$ws.Range("I91").Value = 3000
$ws.Range("J91").Value = 5999.5
$ws.Range("K91").Value = 9000
$ws.Range("L91").Value = 17998.5
$ws.Range("M91").Value = -7518
$ws.Range("N91").Value = -20962.5
$ws.Range("H131").Value = 6956260.5
$ws.Range("J131").Value = 12844.97
$ws.Range("L131").Value = 38534.91
$ws.Range("N131").Value = -48614.91

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 992.3889
$ws.Range("J97").Value = 2028
$ws.Range("L97").Value = 2028
$ws.Range("N97").Value = -3020
$ws.Range("H126").Value = 55640.367
$ws.Range("I126").Value = 3198.9375
$ws.Range("K126").Value = 9596.8125
$ws.Range("M126").Value = -7126.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4324.08
$ws.Range("J7").Value = 5625.8667
$ws.Range("L7").Value = 5625.8667
$ws.Range("N7").Value = -5849.8667
$ws.Range("H46").Value = 1536.3636
$ws.Range("I46").Value = 401
$ws.Range("K46").Value = 401
$ws.Range("M46").Value = -213
$ws.Range("H126").Value = 4324.08
$ws.Range("J126").Value = 5625.8667
$ws.Range("L126").Value = 16877.6001
$ws.Range("N126").Value = -21817.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4099.7617
$ws.Range("I136").Value = 3953.4546
$ws.Range("J136").Value = 4260.7
$ws.Range("K136").Value = 11860.3638
$ws.Range("L136").Value = 12782.1
$ws.Range("M136").Value = -9310.363799999999
$ws.Range("N136").Value = -17882.1
